# MVP and SPL added
#
# Cerro gets an extra SPL Bonus point (row 7), and Maurizio is awarded an
# MVP point (originally row 9). The extra MVP point lifts Maurizio's Total
# above Tito's (originally row 8), so the two players swap places in the
# ranked table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cerro (row 7): SPL Bonus added, derived PointsxG/Total recalculated ---
$ws.Cells.Item(7, 11).Value = 1      # K7  SPL Bonus
$ws.Cells.Item(7, 14).Value = 12.67  # N7  PointsxG
$ws.Cells.Item(7, 15).Value = 76     # O7  Total

# --- Swap the Tito / Maurizio rows (row 8 <-> row 9) via cut/paste so the
#     existing cell formatting & shared-string text cells are preserved ---
$ws.Range("A8:P8").Cut($ws.Range("A37:P37"))
$ws.Range("A9:P9").Cut($ws.Range("A8:P8"))
$ws.Range("A37:P37").Cut($ws.Range("A9:P9"))

# Row 8 now holds Maurizio: MVP point added, rank + derived columns updated
$ws.Cells.Item(8, 2).Value = 7       # B8  Rank
$ws.Cells.Item(8, 10).Value = 1      # J8  MVP
$ws.Cells.Item(8, 14).Value = 10.29  # N8  PointsxG
$ws.Cells.Item(8, 15).Value = 72     # O8  Total
$ws.Cells.Item(8, 16).Value = 4      # P8  Rank Change

# Row 9 now holds Tito: rank-change value updated
$ws.Cells.Item(9, 16).Value = 0      # P9  Rank Change
